# Auto-generated edit script applying scheduled-runner price/profit updates
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H19").Value = 1191
$ws.Range("I19").Value = 718.8571
$ws.Range("J19").Value = 1699.4615
$ws.Range("K19").Value = 718.8571
$ws.Range("L19").Value = 1699.4615
$ws.Range("M19").Value = -543.8571
$ws.Range("N19").Value = -2049.4615

$ws.Range("H38").Value = 245.2
$ws.Range("I38").Value = 42.5
$ws.Range("J38").Value = 1056
$ws.Range("K38").Value = 127.5
$ws.Range("L38").Value = 3168
$ws.Range("M38").Value = 244.5
$ws.Range("N38").Value = -3912

$ws.Range("H74").Value = 3443.4285
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3443.4285
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3443.4285
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5315.4285

$ws.Range("H77").Value = 3443.4285
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3443.4285
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17217.1425
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -26577.1425

$ws.Range("H86").Value = 5265.552
$ws.Range("I86").Value = 5184.381
$ws.Range("J86").Value = 5478.625
$ws.Range("K86").Value = 5184.381
$ws.Range("L86").Value = 5478.625
$ws.Range("M86").Value = -4061.381
$ws.Range("N86").Value = -7724.625

$ws.Range("H89").Value = 5265.552
$ws.Range("I89").Value = 5184.381
$ws.Range("J89").Value = 5478.625
$ws.Range("K89").Value = 25921.905
$ws.Range("L89").Value = 27393.125
$ws.Range("M89").Value = -20305.905
$ws.Range("N89").Value = -38625.125

$ws.Range("H103").Value = 4186.067
$ws.Range("J103").Value = 829.7692
$ws.Range("L103").Value = 2489.3076
$ws.Range("N103").Value = -3661.3076

$ws.Range("H138").Value = 3590744.2
$ws.Range("I138").Value = 1127.6342
$ws.Range("J138").Value = 13402363
$ws.Range("K138").Value = 3382.9026
$ws.Range("L138").Value = 40207089
$ws.Range("M138").Value = 1757.0974
$ws.Range("N138").Value = -40217369

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H32").Value = 5101.523
$ws.Range("I32").Value = 3296.5344
$ws.Range("J32").Value = 20057.143
$ws.Range("K32").Value = 3296.5344
$ws.Range("L32").Value = 20057.143
$ws.Range("M32").Value = -3009.5344
$ws.Range("N32").Value = -20631.143

$ws.Range("H92").Value = 19231.3
$ws.Range("J92").Value = 19231.3
$ws.Range("L92").Value = 19231.3
$ws.Range("N92").Value = -24223.3

$ws.Range("H97").Value = 1508.5555
$ws.Range("I97").Value = 1511
$ws.Range("K97").Value = 1511
$ws.Range("M97").Value = -1015

$ws.Range("H102").Value = 2517.1428
$ws.Range("I102").Value = 2401.8
$ws.Range("J102").Value = 2805.5
$ws.Range("K102").Value = 2401.8
$ws.Range("L102").Value = 2805.5
$ws.Range("M102").Value = -779.8000000000002
$ws.Range("N102").Value = -6049.5

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H9").Value = 25789
$ws.Range("J9").Value = 25789
$ws.Range("L9").Value = 25789
$ws.Range("N9").Value = -26125

$ws.Range("H37").Value = 9121.333
$ws.Range("I37").Value = 2624
$ws.Range("J37").Value = 12370
$ws.Range("K37").Value = 2624
$ws.Range("L37").Value = 12370
$ws.Range("M37").Value = -2487
$ws.Range("N37").Value = -12644

$ws.Range("H44").Value = 22105
$ws.Range("J44").Value = 22105
$ws.Range("L44").Value = 22105
$ws.Range("N44").Value = -23099

$ws.Range("H92").Value = 28940.2
$ws.Range("J92").Value = 28940.2
$ws.Range("L92").Value = 28940.2
$ws.Range("N92").Value = -33932.2

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H23").Value = 17300
$ws.Range("I23").Value = 2500
$ws.Range("K23").Value = 2500
$ws.Range("M23").Value = -2260

$ws.Range("H27").Value = 17300
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2308

$ws.Range("H86").Value = 2717414.2
$ws.Range("I86").Value = 5854572.5
$ws.Range("J86").Value = 8050.273
$ws.Range("K86").Value = 5854572.5
$ws.Range("L86").Value = 8050.273
$ws.Range("M86").Value = -5853449.5
$ws.Range("N86").Value = -10296.273

$ws.Range("H89").Value = 2717414.2
$ws.Range("I89").Value = 5854572.5
$ws.Range("J89").Value = 8050.273
$ws.Range("K89").Value = 29272862.5
$ws.Range("L89").Value = 40251.365
$ws.Range("M89").Value = -29267246.5
$ws.Range("N89").Value = -51483.365

$ws.Range("H141").Value = 29939.111
$ws.Range("J141").Value = 29939.111
$ws.Range("L141").Value = 29939.111
$ws.Range("N141").Value = -40299.111

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H46").Value = 1002442.44
$ws.Range("I46").Value = 633.3333
$ws.Range("J46").Value = 1252894.8
$ws.Range("K46").Value = 1899.9999
$ws.Range("L46").Value = 3758684.4
$ws.Range("M46").Value = -1808.9999
$ws.Range("N46").Value = -3758866.4

$ws.Range("H86").Value = 421.2
$ws.Range("J86").Value = 421.2
$ws.Range("L86").Value = 1263.6
$ws.Range("N86").Value = -3635.6

$ws.Range("H89").Value = 421.2
$ws.Range("J89").Value = 421.2
$ws.Range("L89").Value = 3790.8
$ws.Range("N89").Value = -15646.8

$ws.Range("H113").Value = 866248.06
$ws.Range("I113").Value = 2331481.8
$ws.Range("J113").Value = 428.18182
$ws.Range("K113").Value = 6994445.399999999
$ws.Range("L113").Value = 1284.54546
$ws.Range("M113").Value = -6992275.399999999
$ws.Range("N113").Value = -5624.54546

$ws.Range("H131").Value = 872.25
$ws.Range("I131").Value = 564.25
$ws.Range("J131").Value = 899.0326
$ws.Range("K131").Value = 1692.75
$ws.Range("L131").Value = 2697.0978
$ws.Range("M131").Value = 3347.25
$ws.Range("N131").Value = -12777.0978

$ws.Range("H134").Value = 7497.091
$ws.Range("I134").Value = 4847
$ws.Range("J134").Value = 9705.5
$ws.Range("K134").Value = 14541
$ws.Range("L134").Value = 29116.5
$ws.Range("M134").Value = -9471
$ws.Range("N134").Value = -39256.5

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H57").Value = 17270.066
$ws.Range("I57").Value = 15055
$ws.Range("J57").Value = 17428.285
$ws.Range("K57").Value = 15055
$ws.Range("L57").Value = 17428.285
$ws.Range("M57").Value = -14235
$ws.Range("N57").Value = -19068.285

$ws.Range("H102").Value = 9260333
$ws.Range("I102").Value = 27778068
$ws.Range("J102").Value = 1465
$ws.Range("K102").Value = 27778068
$ws.Range("L102").Value = 1465
$ws.Range("M102").Value = -27776446
$ws.Range("N102").Value = -4709

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H31").Value = 2558.1428
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 2901.1667
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 2901.1667
$ws.Range("M31").Value = -252
$ws.Range("N31").Value = -3397.1667

$ws.Range("H68").Value = 7573.737
$ws.Range("I68").Value = 11817.9
$ws.Range("J68").Value = 2858
$ws.Range("K68").Value = 11817.9
$ws.Range("L68").Value = 2858
$ws.Range("M68").Value = -11068.9
$ws.Range("N68").Value = -4356

$ws.Range("H71").Value = 7573.737
$ws.Range("I71").Value = 11817.9
$ws.Range("J71").Value = 2858
$ws.Range("K71").Value = 59089.5
$ws.Range("L71").Value = 14290
$ws.Range("M71").Value = -55345.5
$ws.Range("N71").Value = -21778

$ws.Range("H112").Value = 31858.54
$ws.Range("J112").Value = 31858.54
$ws.Range("L112").Value = 31858.54
$ws.Range("N112").Value = -34812.54

$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -22350

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H104").Value = 26326.666
$ws.Range("J104").Value = 26326.666
$ws.Range("L104").Value = 26326.666
$ws.Range("N104").Value = -33314.666

$ws.Range("H132").Value = 3388.5483
$ws.Range("I132").Value = 3890.875
$ws.Range("J132").Value = 2852.7334
$ws.Range("K132").Value = 11672.625
$ws.Range("L132").Value = 8558.2002
$ws.Range("M132").Value = -9142.625
$ws.Range("N132").Value = -13618.2002
